$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logBook")

# --- Row 44: new log entry (22:45 - 23:30, 3rd July) ---
# copy formatting from the row above so the new row matches the existing style
$ws.Range("A43:G43").Copy()
$ws.Range("A44:G44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(44, 1).Value = 43
$ws.Cells.Item(44, 2).Value = 44745
$ws.Cells.Item(44, 3).Value = 0.94791666666666663
$ws.Cells.Item(44, 4).Value = 0.97916666666666663
$ws.Range("E44").Formula = "=D44-C44"
$ws.Cells.Item(44, 6).Value = "Code"
$ws.Cells.Item(44, 7).Value = "1. test videos converted to 640_360 resolution"

# --- Row 45: new log entry (00:00 - 2:00am, 4th July) ---
$ws.Range("A43:G43").Copy()
$ws.Range("A45:G45").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = 44745
$ws.Cells.Item(45, 3).Value = 0
$ws.Cells.Item(45, 4).Value = 0.083333333333333329
$ws.Range("E45").Formula = "=D45-C45"
$ws.Cells.Item(45, 6).Value = "Code"
$ws.Cells.Item(45, 7).Value = "1. images and labels converted to 640_480 resolution`n2. deeplabv3+ r50 model train for 12ep on 360 640 dataset`n3. PSPNet r50 model train for 12ep on 360 640 dataset"
$ws.Rows.Item(45).RowHeight = 45

# update selection to match the final state
$ws.Range("E49").Select()

$wb.Save()
